# Updated symbol list on Wed Jan 18 19:47:48 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns for the crypto symbol list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'294.09"
$ws.Range("E2").Value = "'-2.95%"
$ws.Range("E3").Value = "'-3.16%"
$ws.Range("D4").Value = "'4.900"
$ws.Range("E4").Value = "'-1.79%"
$ws.Range("D5").Value = "'0.07315"
$ws.Range("E5").Value = "'-7.55%"
$ws.Range("D6").Value = "'1.813"
$ws.Range("E6").Value = "'-13.81%"
$ws.Range("D7").Value = "'7.692"
$ws.Range("E7").Value = "'-2.07%"
$ws.Range("D8").Value = "'3.761"
$ws.Range("E8").Value = "'-0.89%"
$ws.Range("D9").Value = "'0.9069"
$ws.Range("D10").Value = "'0.1665"
$ws.Range("E10").Value = "'-5.04%"
$ws.Range("D11").Value = "'0.07490"
$ws.Range("E11").Value = "'-6.85%"
$ws.Range("D12").Value = "'0.08104"
$ws.Range("E12").Value = "'-8.09%"
$ws.Range("D13").Value = "'0.02985"
$ws.Range("D15").Value = "'0.001492"
$ws.Range("E15").Value = "'-3.27%"
$ws.Range("D16").Value = "'0.005683"
$ws.Range("E16").Value = "'-4.11%"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.17%"
$ws.Range("E18").Value = "'-7.50%"
$ws.Range("E19").Value = "'-0.40%"
$ws.Range("D20").Value = "'0.1307"
$ws.Range("E20").Value = "'1.32%"
$ws.Range("D21").Value = "'4.324"
$ws.Range("E21").Value = "'4.01%"
$ws.Range("D23").Value = "'0.04476"
$ws.Range("E23").Value = "'-2.76%"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").Value = "'-0.90%"
$ws.Range("D25").Value = "'0.004047"
$ws.Range("E25").Value = "'-10.23%"
$ws.Range("E26").Value = "'0.14%"
$ws.Range("D39").Value = "'0.01653"
$ws.Range("E39").Value = "'-4.61%"
$ws.Range("D40").Value = "'0.04403"
$ws.Range("E40").Value = "'-10.60%"
$ws.Range("D41").Value = "'0.007427"
$ws.Range("E41").Value = "'1.07%"
$ws.Range("E42").Value = "'-3.33%"
$ws.Range("E43").Value = "'-9.83%"
$ws.Range("D44").Value = "'0.01121"
$ws.Range("E44").Value = "'0.95%"
$ws.Range("D45").Value = "'0.00005979"
$ws.Range("E45").Value = "'-1.39%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.23%"
$ws.Range("D47").Value = "'2.157"
$ws.Range("E47").Value = "'162.96%"
$ws.Range("D48").Value = "'0.002434"
$ws.Range("E48").Value = "'-28.37%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.23%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.23%"
